# Add a new "PF/1.0.5" row to the meta-sheet, mirroring the existing
# header/value rows (dev2/sit2/uat2/prod, PF/1.0.0 x4).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "PF/1.0.5"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"

# New row keeps the default/"Normal" cell style (no inherited column
# formatting), same as the other data written straight into row 3.
$ws.Range("A3:D3").Style = "Normal"
